$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64
$ws.Cells.Item($row, 1).Value = "'2025-04-30"
$ws.Cells.Item($row, 2).Value = "substance active"
$ws.Cells.Item($row, 3).Value = 82
$ws.Cells.Item($row, 4).Value = 1
